# Corrected display errors in convergence table
# Rows 5-9 of the convergence synthesis table had several cells showing the
# wrong computed/labelled values. Row 10 is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("F5").Value = "Regular_RightTriangles"
$ws.Range("G5").Value = 0.0212
$ws.Range("H5").Value = "Triangles"
$ws.Range("I5").Value = "Green"
$ws.Range("J5").Value = 17.547

# Row 6
$ws.Range("A6").Value = 6
$ws.Range("F6").Value = "Unstructured_triangles"
$ws.Range("G6").Value = 0.6138
$ws.Range("H6").Value = "Triangles"
$ws.Range("I6").Value = "Green"
$ws.Range("J6").Value = 2.386

# Row 7
$ws.Range("A7").Value = 8
$ws.Range("F7").Value = "Structured_hexagons"
$ws.Range("G7").Value = 1.9416
$ws.Range("H7").Value = "Hexagons"
$ws.Range("I7").Value = "Green"
$ws.Range("J7").Value = 2.38

# Row 8
$ws.Range("A8").Value = 4
$ws.Range("F8").Value = "Equilateral_triangles"
$ws.Range("G8").Value = 1.9766
$ws.Range("H8").Value = "Triangles"
$ws.Range("I8").Value = "Green"
$ws.Range("J8").Value = 5.698

# Row 9
$ws.Range("A9").Value = 0
$ws.Range("F9").Value = "RegularSquares"
$ws.Range("G9").Value = 2.0039
$ws.Range("H9").Value = "Squares"
$ws.Range("I9").Value = "Green"
$ws.Range("J9").Value = 10.755
